# Swap the presentation's theme color palette: the deck currently uses the
# "Integral" theme (green/yellow/teal accents) on its slide master; the
# target state uses the plain default "Office Theme" palette instead (the
# palette that used to live in the notes-master-only theme part). The
# font scheme and format scheme are identical between the two themes, so
# only the 12 theme colours (clrScheme) need to change.
#
# MsoThemeColorSchemeIndex order exposed via ThemeColorScheme.Colors(i):
#   1 dk1   2 lt1   3 dk2   4 lt2   5 accent1  6 accent2
#   7 accent3  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
#
# .RGB expects a BGR-packed integer (classic VB RGB() macro): r + g*256 + b*65536.

function ConvertTo-BgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the standard Office Theme colours.
$officeThemeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $colors.Colors($i).RGB = ConvertTo-BgrInt $officeThemeHex[$i - 1]
}
